$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns D (P10) and E (P90); this shifts the Amazon forecast
# columns (F:I) left into D:G.
$ws.Range("D1:E1").EntireColumn.Delete()

# Update the Prophet Forecast values (column C) with refreshed figures.
$prophet = @(308,406,393,299,207,160,169,184,172,160,161,174,181,172,164,163,176,176,156,137)
for ($i = 0; $i -lt $prophet.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $prophet[$i]
}
